# The uploaded workbook (VAR_param.xlsx) had two typo-ish values in row 4
# corrected: C4 "1.9231" -> "1.19231" and D4 "1.992" -> "1.1992" (a missing
# leading "1" after the decimal point was inserted back in), while B4 and
# E4 keep their original values. The user's cursor ends up resting on E4
# after making the edits.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = "1.19231"
$ws.Range("D4").Value = "1.1992"

$ws.Range("E4").Select()
